$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = "سپهر"
$ws.Range("B6").Value = "09362172533"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = "۱۴۰۴/۷/۹, ۱۷:۵۵:۳۲"

# Row 7
$ws.Range("A7").Value = "سپهر"
$ws.Range("B7").Value = "09362172533"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = "۱۴۰۴/۷/۹, ۱۷:۵۹:۱۰"
